# Apply betexplorer re-scrape update: reshuffle same-date match rows
# and append the newly scraped Diagoras vs Athens Kallithea fixture.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83: swap in the data that now sorts into this row
$ws.Cells.Item(83, 6).Value = "Karditsa"
$ws.Cells.Item(83, 7).Value = 2
$ws.Cells.Item(83, 8).Value = "Levadiakos"
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 4.29
$ws.Cells.Item(83, 11).Value = "18/11/2023 02:12"
$ws.Cells.Item(83, 12).Value = 5.74
$ws.Cells.Item(83, 13).Value = "19/11/2023 13:50"
$ws.Cells.Item(83, 14).Value = 3.08
$ws.Cells.Item(83, 15).Value = "18/11/2023 02:12"
$ws.Cells.Item(83, 16).Value = 3.31
$ws.Cells.Item(83, 17).Value = "19/11/2023 13:50"
$ws.Cells.Item(83, 18).Value = 1.79
$ws.Cells.Item(83, 19).Value = "18/11/2023 02:12"
$ws.Cells.Item(83, 20).Value = 1.68
$ws.Cells.Item(83, 21).Value = "19/11/2023 13:50"
$ws.Cells.Item(83, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/karditsa-levadiakos/2cncBFYt/"

# Row 84: swap in the data that now sorts into this row
$ws.Cells.Item(84, 6).Value = "Kampaniakos"
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = "Iraklis 1908"
$ws.Cells.Item(84, 9).Value = 3
$ws.Cells.Item(84, 10).Value = 3.87
$ws.Cells.Item(84, 11).Value = "18/11/2023 02:12"
$ws.Cells.Item(84, 12).Value = 4.12
$ws.Cells.Item(84, 13).Value = "19/11/2023 12:04"
$ws.Cells.Item(84, 14).Value = 3.15
$ws.Cells.Item(84, 15).Value = "18/11/2023 02:12"
$ws.Cells.Item(84, 16).Value = 3.11
$ws.Cells.Item(84, 17).Value = "19/11/2023 12:52"
$ws.Cells.Item(84, 18).Value = 1.85
$ws.Cells.Item(84, 19).Value = "18/11/2023 02:12"
$ws.Cells.Item(84, 20).Value = 1.97
$ws.Cells.Item(84, 21).Value = "19/11/2023 12:52"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/kampaniakos-iraklis-fc/6R0RGyBJ/"

# Row 85: swap in the data that now sorts into this row
$ws.Cells.Item(85, 6).Value = "Apollon Pontou"
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = "Aiolikos"
$ws.Cells.Item(85, 9).Value = 1
$ws.Cells.Item(85, 10).Value = 2.77
$ws.Cells.Item(85, 11).Value = "19/11/2023 03:12"
$ws.Cells.Item(85, 12).Value = 2.89
$ws.Cells.Item(85, 13).Value = "19/11/2023 13:41"
$ws.Cells.Item(85, 14).Value = 2.87
$ws.Cells.Item(85, 15).Value = "19/11/2023 03:12"
$ws.Cells.Item(85, 16).Value = 2.85
$ws.Cells.Item(85, 17).Value = "19/11/2023 13:24"
$ws.Cells.Item(85, 18).Value = 2.65
$ws.Cells.Item(85, 19).Value = "19/11/2023 03:12"
$ws.Cells.Item(85, 20).Value = 2.67
$ws.Cells.Item(85, 21).Value = "19/11/2023 13:41"
$ws.Cells.Item(85, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/apollon-pontou-aiolikos-fc/pIaVFeQP/"

# Row 87: swap in the data that now sorts into this row
$ws.Cells.Item(87, 6).Value = "Kalamata"
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = "Diagoras"
$ws.Cells.Item(87, 9).Value = 2
$ws.Cells.Item(87, 10).Value = 1.3
$ws.Cells.Item(87, 11).Value = "18/11/2023 02:12"
$ws.Cells.Item(87, 12).Value = 1.26
$ws.Cells.Item(87, 13).Value = "19/11/2023 13:49"
$ws.Cells.Item(87, 14).Value = 4.57
$ws.Cells.Item(87, 15).Value = "18/11/2023 02:12"
$ws.Cells.Item(87, 16).Value = 5.23
$ws.Cells.Item(87, 17).Value = "19/11/2023 13:56"
$ws.Cells.Item(87, 18).Value = 8.08
$ws.Cells.Item(87, 19).Value = "18/11/2023 02:12"
$ws.Cells.Item(87, 20).Value = 12.07
$ws.Cells.Item(87, 21).Value = "19/11/2023 13:56"
$ws.Cells.Item(87, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/kalamata-diagoras-fc/xxnsO5R7/"

# Row 88: swap in the data that now sorts into this row
$ws.Cells.Item(88, 6).Value = "Ionikos"
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = "Athens Kallithea"
$ws.Cells.Item(88, 9).Value = 1
$ws.Cells.Item(88, 10).Value = 3.12
$ws.Cells.Item(88, 11).Value = "18/11/2023 02:12"
$ws.Cells.Item(88, 12).Value = 2.89
$ws.Cells.Item(88, 13).Value = "19/11/2023 13:29"
$ws.Cells.Item(88, 14).Value = 2.85
$ws.Cells.Item(88, 15).Value = "18/11/2023 02:12"
$ws.Cells.Item(88, 16).Value = 2.91
$ws.Cells.Item(88, 17).Value = "19/11/2023 13:29"
$ws.Cells.Item(88, 18).Value = 2.26
$ws.Cells.Item(88, 19).Value = "18/11/2023 02:12"
$ws.Cells.Item(88, 20).Value = 2.62
$ws.Cells.Item(88, 21).Value = "19/11/2023 13:29"
$ws.Cells.Item(88, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/ionikos-athens-kallithea/v1SypmlR/"

# Row 89: swap in the data that now sorts into this row
$ws.Cells.Item(89, 6).Value = "Giouchtas"
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = "PAE Chania"
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 3.28
$ws.Cells.Item(89, 11).Value = "18/11/2023 02:12"
$ws.Cells.Item(89, 12).Value = 3.93
$ws.Cells.Item(89, 13).Value = "19/11/2023 13:14"
$ws.Cells.Item(89, 14).Value = 2.86
$ws.Cells.Item(89, 15).Value = "18/11/2023 02:12"
$ws.Cells.Item(89, 16).Value = 3.05
$ws.Cells.Item(89, 17).Value = "19/11/2023 12:51"
$ws.Cells.Item(89, 18).Value = 2.17
$ws.Cells.Item(89, 19).Value = "18/11/2023 02:12"
$ws.Cells.Item(89, 20).Value = 2.04
$ws.Cells.Item(89, 21).Value = "19/11/2023 13:14"
$ws.Cells.Item(89, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/giouchtas-pae-chania/4rzPnk48/"

# Row 91: swap in the data that now sorts into this row
$ws.Cells.Item(91, 6).Value = "Iraklis 1908"
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = "AEK Athens FC B"
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 1.5
$ws.Cells.Item(91, 11).Value = "24/11/2023 02:12"
$ws.Cells.Item(91, 12).Value = 1.67
$ws.Cells.Item(91, 13).Value = "25/11/2023 13:45"
$ws.Cells.Item(91, 14).Value = 3.67
$ws.Cells.Item(91, 15).Value = "24/11/2023 02:12"
$ws.Cells.Item(91, 16).Value = 3.53
$ws.Cells.Item(91, 17).Value = "25/11/2023 13:45"
$ws.Cells.Item(91, 18).Value = 5.66
$ws.Cells.Item(91, 19).Value = "24/11/2023 02:12"
$ws.Cells.Item(91, 20).Value = 5.3
$ws.Cells.Item(91, 21).Value = "25/11/2023 13:45"
$ws.Cells.Item(91, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/iraklis-fc-aek/zmvA8DIa/"

# Row 93: swap in the data that now sorts into this row
$ws.Cells.Item(93, 6).Value = "Athens Kallithea"
$ws.Cells.Item(93, 7).Value = 3
$ws.Cells.Item(93, 8).Value = "Tilikratis L."
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 1.16
$ws.Cells.Item(93, 11).Value = "24/11/2023 02:12"
$ws.Cells.Item(93, 12).Value = 1.11
$ws.Cells.Item(93, 13).Value = "25/11/2023 13:09"
$ws.Cells.Item(93, 14).Value = 6.22
$ws.Cells.Item(93, 15).Value = "24/11/2023 02:12"
$ws.Cells.Item(93, 16).Value = 8.29
$ws.Cells.Item(93, 17).Value = "25/11/2023 13:10"
$ws.Cells.Item(93, 18).Value = 12.53
$ws.Cells.Item(93, 19).Value = "24/11/2023 02:12"
$ws.Cells.Item(93, 20).Value = 23.73
$ws.Cells.Item(93, 21).Value = "25/11/2023 13:10"
$ws.Cells.Item(93, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/athens-kallithea-tilikratis-lefkada/QL6KyPe7/"

# Row 95: swap in the data that now sorts into this row
$ws.Cells.Item(95, 6).Value = "Makedonikos"
$ws.Cells.Item(95, 7).Value = 3
$ws.Cells.Item(95, 8).Value = "Apollon Pontou"
$ws.Cells.Item(95, 9).Value = 1
$ws.Cells.Item(95, 10).Value = 1.53
$ws.Cells.Item(95, 11).Value = "25/11/2023 02:12"
$ws.Cells.Item(95, 12).Value = 1.54
$ws.Cells.Item(95, 13).Value = "26/11/2023 13:44"
$ws.Cells.Item(95, 14).Value = 3.51
$ws.Cells.Item(95, 15).Value = "25/11/2023 02:12"
$ws.Cells.Item(95, 16).Value = 3.65
$ws.Cells.Item(95, 17).Value = "26/11/2023 13:59"
$ws.Cells.Item(95, 18).Value = 5.58
$ws.Cells.Item(95, 19).Value = "25/11/2023 02:12"
$ws.Cells.Item(95, 20).Value = 6.81
$ws.Cells.Item(95, 21).Value = "26/11/2023 13:59"
$ws.Cells.Item(95, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/makedonikos-neapolis-apollon-pontou/QyXM5B2I/"

# Row 96: swap in the data that now sorts into this row
$ws.Cells.Item(96, 6).Value = "Ilioupoli"
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = "Giouchtas"
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 1.97
$ws.Cells.Item(96, 11).Value = "25/11/2023 02:12"
$ws.Cells.Item(96, 12).Value = 2.19
$ws.Cells.Item(96, 13).Value = "26/11/2023 13:59"
$ws.Cells.Item(96, 14).Value = 3.02
$ws.Cells.Item(96, 15).Value = "25/11/2023 02:12"
$ws.Cells.Item(96, 16).Value = 3.19
$ws.Cells.Item(96, 17).Value = "26/11/2023 13:59"
$ws.Cells.Item(96, 18).Value = 3.62
$ws.Cells.Item(96, 19).Value = "25/11/2023 02:12"
$ws.Cells.Item(96, 20).Value = 3.32
$ws.Cells.Item(96, 21).Value = "26/11/2023 13:59"
$ws.Cells.Item(96, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/ilioupoli-giouchtas/Wb7Gx5t1/"

# Row 98: swap in the data that now sorts into this row
$ws.Cells.Item(98, 6).Value = "PAE Chania"
$ws.Cells.Item(98, 7).Value = 2
$ws.Cells.Item(98, 8).Value = "PAE Egaleo"
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 1.55
$ws.Cells.Item(98, 11).Value = "25/11/2023 02:12"
$ws.Cells.Item(98, 12).Value = 1.2
$ws.Cells.Item(98, 13).Value = "26/11/2023 13:55"
$ws.Cells.Item(98, 14).Value = 3.48
$ws.Cells.Item(98, 15).Value = "25/11/2023 02:12"
$ws.Cells.Item(98, 16).Value = 6.56
$ws.Cells.Item(98, 17).Value = "26/11/2023 13:58"
$ws.Cells.Item(98, 18).Value = 5.44
$ws.Cells.Item(98, 19).Value = "25/11/2023 02:12"
$ws.Cells.Item(98, 20).Value = 11.76
$ws.Cells.Item(98, 21).Value = "26/11/2023 13:58"
$ws.Cells.Item(98, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/pae-chania-pae-egaleo/zXD7vRBl/"

# Row 99: swap in the data that now sorts into this row
$ws.Cells.Item(99, 6).Value = "Karditsa"
$ws.Cells.Item(99, 7).Value = 2
$ws.Cells.Item(99, 8).Value = "Niki Volos"
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 3.26
$ws.Cells.Item(99, 11).Value = "25/11/2023 02:12"
$ws.Cells.Item(99, 12).Value = 3.72
$ws.Cells.Item(99, 13).Value = "26/11/2023 13:46"
$ws.Cells.Item(99, 14).Value = 2.83
$ws.Cells.Item(99, 15).Value = "25/11/2023 02:12"
$ws.Cells.Item(99, 16).Value = 2.91
$ws.Cells.Item(99, 17).Value = "26/11/2023 13:46"
$ws.Cells.Item(99, 18).Value = 2.2
$ws.Cells.Item(99, 19).Value = "25/11/2023 02:12"
$ws.Cells.Item(99, 20).Value = 2.18
$ws.Cells.Item(99, 21).Value = "26/11/2023 13:46"
$ws.Cells.Item(99, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/karditsa-niki-volos/vuyR4VHO/"

# Row 105: new fixture appended at the end of the sheet
$ws.Range("A104:V104").Copy($ws.Range("A105:V105"))
$ws.Cells.Item(105, 1).Value = 104
$ws.Cells.Item(105, 2).Value = "greece"
$ws.Cells.Item(105, 3).Value = "super-league-2"
$ws.Cells.Item(105, 4).Value = "2023-2024"
$ws.Cells.Item(105, 5).Value = 45262.625
$ws.Cells.Item(105, 6).Value = "Diagoras"
$ws.Cells.Item(105, 7).Value = 2
$ws.Cells.Item(105, 8).Value = "Athens Kallithea"
$ws.Cells.Item(105, 9).Value = 1
$ws.Cells.Item(105, 10).Value = 5.88
$ws.Cells.Item(105, 11).Value = "01/12/2023 03:12"
$ws.Cells.Item(105, 12).Value = 4.71
$ws.Cells.Item(105, 13).Value = "02/12/2023 14:54"
$ws.Cells.Item(105, 14).Value = 3.52
$ws.Cells.Item(105, 15).Value = "01/12/2023 03:12"
$ws.Cells.Item(105, 16).Value = 3.1
$ws.Cells.Item(105, 17).Value = "02/12/2023 14:54"
$ws.Cells.Item(105, 18).Value = 1.53
$ws.Cells.Item(105, 19).Value = "01/12/2023 03:12"
$ws.Cells.Item(105, 20).Value = 1.86
$ws.Cells.Item(105, 21).Value = "02/12/2023 14:54"
$ws.Cells.Item(105, 22).Value = "https://www.betexplorer.com/football/greece/super-league-2/diagoras-fc-athens-kallithea/nayJKsNB/"
